$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply corrected cell values per row (restoring the correct pairing of
# item code / rate / qty / value figures that had been shuffled between rows).

$ws.Range("B256").Value = 48719
$ws.Range("E256").Value = 353.35
$ws.Range("F256").Value = -81
$ws.Range("G256").Value = -23955.75

$ws.Range("B257").Value = 64979
$ws.Range("E257").Value = 314.41
$ws.Range("F257").Value = 82
$ws.Range("G257").Value = 24251.5

$ws.Range("B308").Value = 63565
$ws.Range("D308").Value = 102.71
$ws.Range("E308").Value = 109.19
$ws.Range("F308").Value = 60
$ws.Range("G308").Value = 6162.6

$ws.Range("B309").Value = 57077
$ws.Range("D309").Value = 93.08
$ws.Range("E309").Value = 111.2
$ws.Range("F309").Value = 1
$ws.Range("G309").Value = 93.08

$ws.Range("B310").Value = 61610
$ws.Range("E310").Value = 122.71
$ws.Range("F310").Value = -58
$ws.Range("G310").Value = -5957.18

$ws.Range("B342").Value = 57802
$ws.Range("E342").Value = 162.71
$ws.Range("F342").Value = -79
$ws.Range("G342").Value = -11334.92

$ws.Range("B343").Value = 63531
$ws.Range("F343").Value = 80
$ws.Range("G343").Value = 11478.4

$ws.Range("B344").Value = 63571
$ws.Range("E344").Value = 152.53
$ws.Range("F344").Value = 29
$ws.Range("G344").Value = 4160.92

$ws.Range("B347").Value = 63510
$ws.Range("E347").Value = 50.66
$ws.Range("F347").Value = 167
$ws.Range("G347").Value = 7955.88

$ws.Range("B348").Value = 55356
$ws.Range("E348").Value = 54.04
$ws.Range("F348").Value = -158
$ws.Range("G348").Value = -7527.12

$ws.Range("B367").Value = 61605
$ws.Range("E367").Value = 133.78
$ws.Range("F367").Value = -13
$ws.Range("G367").Value = -1455.48

$ws.Range("B368").Value = 63563
$ws.Range("E368").Value = 119.04
$ws.Range("F368").Value = 15
$ws.Range("G368").Value = 1679.4

$ws.Range("B374").Value = 63560
$ws.Range("E374").Value = 134.87
$ws.Range("F374").Value = 104
$ws.Range("G374").Value = 13193.44

$ws.Range("B375").Value = 60325
$ws.Range("E375").Value = 151.57
$ws.Range("F375").Value = -102
$ws.Range("G375").Value = -12939.72

$ws.Range("B381").Value = 57817
$ws.Range("F381").Value = 3
$ws.Range("G381").Value = 239.43

$ws.Range("B382").Value = 62865
$ws.Range("F382").Value = 151
$ws.Range("G382").Value = 12051.31

$ws.Range("B449").Value = 31930
$ws.Range("E449").Value = 26.8
$ws.Range("F449").Value = -62
$ws.Range("G449").Value = -1390.04

$ws.Range("B450").Value = 63681
$ws.Range("E450").Value = 23.84
$ws.Range("F450").Value = 65
$ws.Range("G450").Value = 1457.3

$ws.Range("B528").Value = 47097
$ws.Range("D528").Value = 112.28
$ws.Range("E528").Value = 134.16
$ws.Range("F528").Value = 15
$ws.Range("G528").Value = 1684.2

$ws.Range("B529").Value = 58047
$ws.Range("D529").Value = 105.54
$ws.Range("E529").Value = 126.1
$ws.Range("F529").Value = 54
$ws.Range("G529").Value = 5699.16

$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9

$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68

$ws.Range("B593").Value = 45702
$ws.Range("E593").Value = 31.43
$ws.Range("F593").Value = -215
$ws.Range("G593").Value = -5654.5

$ws.Range("B594").Value = 64919
$ws.Range("E594").Value = 27.97
$ws.Range("F594").Value = 224
$ws.Range("G594").Value = 5891.2

$ws.Range("B596").Value = 53595
$ws.Range("E596").Value = 17.61
$ws.Range("F596").Value = -335
$ws.Range("G596").Value = -4934.55

$ws.Range("B597").Value = 65067
$ws.Range("E597").Value = 15.65
$ws.Range("F597").Value = 338
$ws.Range("G597").Value = 4978.74

$ws.Range("B679").Value = 64810
$ws.Range("E679").Value = 291.22
$ws.Range("F679").Value = 7
$ws.Range("G679").Value = 1917.44

$ws.Range("B680").Value = 53319
$ws.Range("E680").Value = 310.64
$ws.Range("F680").Value = -6
$ws.Range("G680").Value = -1643.52

$ws.Range("B701").Value = 60025
$ws.Range("E701").Value = 37.22
$ws.Range("F701").Value = -98
$ws.Range("G701").Value = -3217.34

$ws.Range("B702").Value = 64833
$ws.Range("E702").Value = 34.9
$ws.Range("F702").Value = 99
$ws.Range("G702").Value = 3250.17

$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79

$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11

$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53

$ws.Range("B865").Value = 65079
$ws.Range("E865").Value = 43.44
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27
